# Generate Report for Handoff
# Adds two new handed-off files (2bdc51fe-... and 61a15d7d-...) as new rows
# to the Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" (table3, displayName "Overview")
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)

$loOverview.ListRows.Add() | Out-Null
$wsOverview.Range("A4").Value = "2bdc51fe-c985-4dcb-995d-f1ac8ee5abb8.md"
$wsOverview.Range("B4").Value = "e2e\2bdc51fe-c985-4dcb-995d-f1ac8ee5abb8.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2016-09-08 04:54:52"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2bdc51fec985c4db5995df1ac8ee5abb8000001/e2e/2bdc51fe-c985-4dcb-995d-f1ac8ee5abb8.md", "", "", "e2e\2bdc51fe-c985-4dcb-995d-f1ac8ee5abb8.md") | Out-Null

$loOverview.ListRows.Add() | Out-Null
$wsOverview.Range("A5").Value = "61a15d7d-b6d5-4da7-b456-f9204bdc3269.md"
$wsOverview.Range("B5").Value = "e2e\61a15d7d-b6d5-4da7-b456-f9204bdc3269.md"
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-09-08 04:54:52"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/61a15d7db6d54da7b456f9204bdc326900000002/e2e/61a15d7d-b6d5-4da7-b456-f9204bdc3269.md", "", "", "e2e\61a15d7d-b6d5-4da7-b456-f9204bdc3269.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (table1)
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)

$loZhCn.ListRows.Add() | Out-Null
$wsZhCn.Range("A4").Value = "2bdc51fe-c985-4dcb-995d-f1ac8ee5abb8.md"
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "False"
$wsZhCn.Range("G4").Value = "2bdc51fe-c985-4dcb-995d-f1ac8ee5abb8.2f0fbf8f307deebf9610aa891686c52304d5e51e.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-09-08 04:54:47"
$wsZhCn.Range("K4").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("M4").Value = "True"
$wsZhCn.Range("O4").Value = "False"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/2bdc51fec985c4db5995df1ac8ee5abb8000003/e2e/2bdc51fe-c985-4dcb-995d-f1ac8ee5abb8.md", "", "", "2bdc51fe-c985-4dcb-995d-f1ac8ee5abb8.md") | Out-Null

$loZhCn.ListRows.Add() | Out-Null
$wsZhCn.Range("A5").Value = "61a15d7d-b6d5-4da7-b456-f9204bdc3269.md"
$wsZhCn.Range("B5").Value = ".md"
$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("D5").Value = "e2e"
$wsZhCn.Range("E5").Value = "ht"
$wsZhCn.Range("F5").Value = "False"
$wsZhCn.Range("G5").Value = "61a15d7d-b6d5-4da7-b456-f9204bdc3269.1d693acdd3bdb2b02ece3b280acdefbee359af93.zh-cn.xlf"
$wsZhCn.Range("H5").Value = "2016-09-08 04:54:47"
$wsZhCn.Range("K5").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("M5").Value = "True"
$wsZhCn.Range("O5").Value = "False"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/61a15d7db6d54da7b456f9204bdc326900000004/e2e/61a15d7d-b6d5-4da7-b456-f9204bdc3269.md", "", "", "61a15d7d-b6d5-4da7-b456-f9204bdc3269.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de" (table2)
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)

$loDeDe.ListRows.Add() | Out-Null
$wsDeDe.Range("A4").Value = "2bdc51fe-c985-4dcb-995d-f1ac8ee5abb8.md"
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "False"
$wsDeDe.Range("G4").Value = "2bdc51fe-c985-4dcb-995d-f1ac8ee5abb8.2f0fbf8f307deebf9610aa891686c52304d5e51e.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-09-08 04:54:52"
$wsDeDe.Range("K4").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("M4").Value = "True"
$wsDeDe.Range("O4").Value = "False"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/2bdc51fec985c4db5995df1ac8ee5abb8000005/e2e/2bdc51fe-c985-4dcb-995d-f1ac8ee5abb8.md", "", "", "2bdc51fe-c985-4dcb-995d-f1ac8ee5abb8.md") | Out-Null

$loDeDe.ListRows.Add() | Out-Null
$wsDeDe.Range("A5").Value = "61a15d7d-b6d5-4da7-b456-f9204bdc3269.md"
$wsDeDe.Range("B5").Value = ".md"
$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("D5").Value = "e2e"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("F5").Value = "False"
$wsDeDe.Range("G5").Value = "61a15d7d-b6d5-4da7-b456-f9204bdc3269.1d693acdd3bdb2b02ece3b280acdefbee359af93.de-de.xlf"
$wsDeDe.Range("H5").Value = "2016-09-08 04:54:52"
$wsDeDe.Range("K5").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("M5").Value = "True"
$wsDeDe.Range("O5").Value = "False"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/61a15d7db6d54da7b456f9204bdc326900000006/e2e/61a15d7d-b6d5-4da7-b456-f9204bdc3269.md", "", "", "61a15d7d-b6d5-4da7-b456-f9204bdc3269.md") | Out-Null
